$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Fill in row 12 with the new user record
$ws.Range("A12").Value = "U1412"
$ws.Range("B12").Value = "آقای عطایی"
$ws.Range("C12").Value = "ataee"
$ws.Range("D12").Value = 123456
$ws.Range("E12").Value = "user"

# Move the active selection to D13
$ws.Range("D13").Select()
